$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Expand the "would check ... STAT register" sentence into the new,
# longer explanation, add the new sentences about positive/negative branch
# instructions and taking a branch, and start a new paragraph (with an empty
# paragraph in between) that begins "The logic to change the program
# counter for branching was given in a file called ".
# ---------------------------------------------------------------------------
$found1 = $d.Content.Find.Execute(
  "would check using the STAT register. If a branch was unconditional we would have to take it but wouldn" + [char]8217 + "t have to compare the STAT register. To add the branches we were given a file called ",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "would checked using the STAT register. If we were in a positive branch instruction (BRA or BRR), we would take the branch only if the zero bit of the STAT register is set to 1.  On the other hand, for negative branch instructions (BNE or BNR), we take the branch only if the zero bit of the STAT register is set to 0.  Taking a branch involved changing the program counter to a new value relatively (BRR or BNR) or absolutely (BNE or BNR) depending on the branch type.  If a branch is not taken, the program counter does not change until the beginning of the next instruction.[PARABREAK]The logic to change the program counter for branching was given in a file called ",
  2)

# Turn the marker into the two actual paragraph breaks (one blank paragraph)
# as a separate, tightly-scoped replace so it doesn't disturb the run
# boundaries of the text around it (in particular the "Br.v" run that
# immediately follows).
$found1b = $d.Content.Find.Execute("[PARABREAK]", $true, $false, $false, $false, $false, $true, 1, $false, "^p^p", 2)

# ---------------------------------------------------------------------------
# Step 2: After "Br.v" (1st occurrence), expand "that was used to calculate
# ... the program counter." and add the "Br.v took as input br_sel ..."
# sentence before "The program counter (PC) was held in ".
# ---------------------------------------------------------------------------
$found2 = $d.Content.Find.Execute(
  "that was used to calculate the location of the branch and would then send it to the program counter. The program counter (PC) was held in ",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "that was used to calculate the location of the branch and would then send it to the program counter.  Br.v took as input br_sel from our ctrl.v file which determined if the program counter would be modified in a relative or absolute calculation.  The program counter (PC) was held in ",
  2)

# ---------------------------------------------------------------------------
# Step 3: After "PC.v", fix "another file we were give" -> "given", and
# add the "PC.v controls if a branch is taken ..." sentence before
# "When we started this section of the project ...".
# ---------------------------------------------------------------------------
$found3 = $d.Content.Find.Execute(
  ", another file we were give, it would control if a branch was taken or the PC was only incremented. When we started this section of the project we added all the new files we were given to the ",
  $true, $false, $false, $false, $false, $true, 1, $false,
  ", another file we were given.  PC.v controls if a branch is taken or the PC was only incremented based on the input pc_sel whose value is set in our ctrl.v file. When we started this section of the project we added all the new files we were given to the ",
  2)

# ---------------------------------------------------------------------------
# Step 4: After "SISC.v file ... via wires.", add the new sentence about the
# imem.data file.
# ---------------------------------------------------------------------------
$found4 = $d.Content.Find.Execute(
  " file as new module instantiations and then connected the inputs and outputs via wires. ",
  $true, $false, $false, $false, $false, $true, 1, $false,
  " file as new module instantiations and then connected the inputs and outputs via wires.  Among the new files was the imem.data file which provided instructions to test our implementation.  These instructions were run,",
  2)

# ---------------------------------------------------------------------------
# Step 5: The final run in the document (right after the _GoBack bookmark)
# was a single space; replace it with the closing sentence about imem.data
# outputs. We target it directly via a Range at the very end of the
# document content (rather than Find) because Find does not reliably keep
# its match isolated right after the bookmark boundary.
# ---------------------------------------------------------------------------
$tailRange = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$tailRange.Text = " and the outputs matched those expected as indicated in the comments from the imem.data file."
